# "cierre de 7 sept 2021" — close out the current voucher ("VALES DE
# INSENTIVOS") sheet: new amount, new amount-in-words, new concept, new
# recipient name. The TODAY() date stamp in A11 is volatile and recalculates
# on its own; we just leave the cursor where the clerk would after filling
# the form in (the merged signature-date row, A10:D10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Amount paid (currency-formatted) and its spelled-out equivalent.
$ws.Range("D1").Value = 5000
$ws.Range("A2").Value = "CINCO MIL     PESOS 00/100 M.N."

# What the payment is for.
$ws.Range("A4").Value = "PAGO DE INCENTIVO DEL MES DE AGOSTO 2021"

# Who is receiving it.
$ws.Range("C8").Value = "PABLO BAEZ"

# Leave the selection on the signature/date row, like after tabbing through
# the form.
$ws.Range("A10:D10").Select() | Out-Null
